$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("B9").Value = "8:41 AM"
$ws.Range("C9").Value = 11729.405
$ws.Range("D9").Value = 11759.405
$ws.Range("F9").Value = 21.885

# Row 10
$ws.Range("A10").Value = 396
$ws.Range("B10").Value = "8:43 AM"
$ws.Range("C10").Value = 11862.452
$ws.Range("D10").Value = 11892.452
$ws.Range("F10").Value = 11.024

# Row 11
$ws.Range("A11").Value = 399
$ws.Range("C11").Value = 11968.34
$ws.Range("D11").Value = 11998.34
$ws.Range("F11").Value = 13.115

# Row 17
$ws.Range("A17").Value = 436
$ws.Range("C17").Value = 13056.61
$ws.Range("D17").Value = 13086.61
$ws.Range("E17").Value = 3.23
$ws.Range("F17").Value = 13.92

# Row 19
$ws.Range("A19").Value = 444
$ws.Range("B19").Value = "9:07 AM"
$ws.Range("C19").Value = 13309.11
$ws.Range("D19").Value = 13339.11
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 67.56
$ws.Range("G19").Value = 0

# Row 25
$ws.Range("A25").Value = 467
$ws.Range("B25").Value = "9:18 AM"
$ws.Range("C25").Value = 14001.8
$ws.Range("D25").Value = 14031.8
$ws.Range("F25").Value = 34.39

# Row 26
$ws.Range("A26").Value = 474
$ws.Range("B26").Value = "9:22 AM"
$ws.Range("C26").Value = 14204.495
$ws.Range("D26").Value = 14234.495
$ws.Range("F26").Value = 40.96

# Row 27
$ws.Range("A27").Value = 480
$ws.Range("B27").Value = "9:25 AM"
$ws.Range("C27").Value = 14376.725
$ws.Range("D27").Value = 14406.725
$ws.Range("F27").Value = 15.025

# Row 28
$ws.Range("A28").Value = 484
$ws.Range("B28").Value = "9:27 AM"
$ws.Range("C28").Value = 14495.41
$ws.Range("D28").Value = 14525.41
$ws.Range("F28").Value = 65.56999999999999

# Row 30
$ws.Range("A30").Value = 497
$ws.Range("C30").Value = 14894.885
$ws.Range("D30").Value = 14924.885
$ws.Range("F30").Value = 30.35

# Row 31
$ws.Range("A31").Value = 501
$ws.Range("B31").Value = "9:35 AM"
$ws.Range("C31").Value = 15000.94
$ws.Range("D31").Value = 15030.94
$ws.Range("F31").Value = 51.39

# Row 32
$ws.Range("A32").Value = 504
$ws.Range("B32").Value = "9:37 AM"
$ws.Range("C32").Value = 15106.36
$ws.Range("D32").Value = 15136.36
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = 52.755
$ws.Range("G32").Value = 0

# Row 33
$ws.Range("A33").Value = 510
$ws.Range("B33").Value = "9:40 AM"
$ws.Range("C33").Value = 15282.67
$ws.Range("D33").Value = 15312.67
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 24.24

# Row 34
$ws.Range("A34").Value = 518
$ws.Range("B34").Value = "9:44 AM"
$ws.Range("C34").Value = 15517.59
$ws.Range("D34").Value = 15547.59
$ws.Range("E34").Value = 20.516667
$ws.Range("F34").Value = 0.58
$ws.Range("G34").Value = 0

# Row 35
$ws.Range("A35").Value = 568
$ws.Range("B35").Value = "10:09 A"
$ws.Range("C35").Value = 17019.785
$ws.Range("D35").Value = 17049.785
$ws.Range("E35").Value = 8.845000000000001
$ws.Range("F35").Value = 1.025
$ws.Range("G35").Value = 1

# Row 36
$ws.Range("A36").Value = 576
$ws.Range("B36").Value = "10:13 A"
$ws.Range("C36").Value = 17250.1
$ws.Range("D36").Value = 17280.1
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 3.83
$ws.Range("G36").Value = 0

# Row 37
$ws.Range("A37").Value = 581
$ws.Range("B37").Value = "10:15 A"
$ws.Range("C37").Value = 17413.61
$ws.Range("D37").Value = 17443.61
$ws.Range("F37").Value = 676.14

# Row 38
$ws.Range("A38").Value = 983
$ws.Range("B38").Value = "1:36 PM"
$ws.Range("C38").Value = 29464.66
$ws.Range("D38").Value = 29494.66
$ws.Range("E38").Value = 0.803333
$ws.Range("F38").Value = 0.78

# Row 39
$ws.Range("A39").Value = 999
$ws.Range("B39").Value = "1:44 PM"
$ws.Range("C39").Value = 29954.343333
$ws.Range("D39").Value = 29984.343333
$ws.Range("E39").Value = 3.856667
$ws.Range("F39").Value = 0.826667
$ws.Range("G39").Value = 0.666667

# Row 40
$ws.Range("A40").Value = 1035
$ws.Range("B40").Value = "2:02 PM"
$ws.Range("C40").Value = 31029.245
$ws.Range("D40").Value = 31059.245
$ws.Range("E40").Value = 17.475
$ws.Range("F40").Value = 0.545
$ws.Range("G40").Value = 0

# Row 44
$ws.Range("A44").Value = 1389
$ws.Range("B44").Value = "4:59 PM"
$ws.Range("C44").Value = 41640.85
$ws.Range("D44").Value = 41670.85
$ws.Range("F44").Value = 3.15

# Row 45
$ws.Range("B45").Value = "5:04 PM"
$ws.Range("C45").Value = 41905.92
$ws.Range("D45").Value = 41935.92
$ws.Range("E45").Value = 7.78
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 1

